# Revert "Default NH3 for more states added"
# - Set the Generic NH3 Emission Rate (column G) back to 0 for the rows that
#   had been given a non-zero default value.
# - Remove the leftover "_xlnm._FilterDatabase" defined name.
# - Restore the sheet view to its plain state (no stale scroll position /
#   selection left over from the filtering session).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rows (in column G) whose NH3 emission rate default needs to go back to 0.
$rows = @(2,4,6,8,15,17,18,20,22,23,25,27,28,30,32,33,35,37,38,40,42,43,45,47,48,50,52,53,55,57,58,60,62,63,65,67,68,70,71,73,75,76,78,80,82,84,86,88,89,91,93,94,96,98,99,102,103,104,106,107,109,111,112,115,116,118,120,121,123,125,126,129,131,132,134,135,137,139,141,142,144,146,147,149,150,152,153,155,157)

foreach ($r in $rows) {
    $ws.Cells.Item($r, 7).Value = 0
}

# Drop the hidden AutoFilter defined name that is no longer needed.
try {
    $fd = $wb.Names.Item("_xlnm._FilterDatabase")
    $fd.Delete()
} catch {
    $names = $wb.Names
    $cnt = $names.Count()
    for ($i = $cnt; $i -ge 1; $i--) {
        $n = $names.Item($i)
        if ($n.Name() -like "*FilterDatabase*") {
            $n.Delete()
        }
    }
}

# Reset the view: scroll back to the top-left corner and select A1 so no
# stale "topLeftCell"/"selection" pointing at G1/H1 remains.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("A1").Select()
